# Fruta / hortaliza, semanal
# Permute the D,L,M,N,O,P,Q,S,T values across rows 2-20 (row 9 stays put).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row r gets the data currently stored in row mapping[r]
$map = @{
    2  = 11
    3  = 13
    4  = 14
    5  = 15
    6  = 20
    7  = 5
    8  = 2
    9  = 9
    10 = 17
    11 = 8
    12 = 10
    13 = 7
    14 = 12
    15 = 16
    16 = 6
    17 = 19
    18 = 3
    19 = 18
    20 = 4
}

$cols = @("D", "L", "M", "N", "O", "P", "Q", "S", "T")

# Snapshot the original values for the columns that move, per row.
$snapshot = @{}
for ($r = 2; $r -le 20; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Write back using the permutation map (source row -> destination row).
foreach ($r in 2..20) {
    $src = $map[$r]
    $srcVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $srcVals[$c]
    }
}
